$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '67.343.69'
$ws.Range('E2').Value = '  -1.45%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.506.13'
$ws.Range('E3').Value = '  -3.58%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '200.16'
$ws.Range('E5').Value = '  +1.57%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '552.48'
$ws.Range('E6').Value = '  -4.45%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '3.496.77'
$ws.Range('E7').Value = '  -3.67%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.607'
$ws.Range('E8').Value = '  -2.05%  '
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.654'
$ws.Range('E10').Value = '  -3.60%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '62.20'
$ws.Range('E11').Value = '  +10.50%  '
$ws.Range('E12').Value = '  -6.95%  '
$ws.Range('E13').Value = '  -7.30%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '9.81'
$ws.Range('E14').Value = '  -2.54%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '4.067.72'
$ws.Range('E15').Value = '  -3.60%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.507.83'
$ws.Range('E16').Value = '  -3.75%  '
$ws.Range('E17').Value = '  -1.98%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '18.44'
$ws.Range('E18').Value = '  -1.10%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '67.096.57'
$ws.Range('E19').Value = '  -1.79%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '11.78'
$ws.Range('E20').Value = '  -6.19%  '
$ws.Range('E21').Value = '  -5.61%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '391.46'
$ws.Range('E22').Value = '  -2.77%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '12.17'
$ws.Range('E23').Value = '  -6.56%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '3.99'
$ws.Range('E24').Value = '  -5.89%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '83.05'
$ws.Range('E25').Value = '  -3.42%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '3.94'
$ws.Range('E26').Value = '  +2.21%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.83'
$ws.Range('E29').Value = '  -3.63%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '30.95'
$ws.Range('E30').Value = '  -2.45%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '692.48'
$ws.Range('E31').Value = '  +0.41%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '7.01'
$ws.Range('E32').Value = '  -13.25%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '11.70'
$ws.Range('E33').Value = '  -4.34%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '63.77'
$ws.Range('E34').Value = '  -1.50%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.111'
$ws.Range('E35').Value = '  -6.13%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '38.60'
$ws.Range('E36').Value = '  -9.68%  '
$ws.Range('E37').Value = '  +0.20%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.395'
$ws.Range('E38').Value = '  -6.38%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0₃0677'
$ws.Range('E43').Value = '  -13.76%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.56'
$ws.Range('E44').Value = '  -10.46%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.76'
$ws.Range('E45').Value = '  +5.33%  '
$ws.Range('E46').Value = '  -4.20%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.66'
$ws.Range('E47').Value = '  -10.14%  '
$ws.Range('E48').Value = '  -3.65%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '138.36'
$ws.Range('E49').Value = '  -2.68%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '8.21'
$ws.Range('E50').Value = '  -7.69%  '
$ws.Range('E51').Value = '  -7.16%  '

# Row swaps (coin identity + price + volume moved between adjacent rows)
$ws.Range('B27').Value = 'ImmutableX'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.81'
$ws.Range('E27').Value = '  -5.01%  '

$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '12.19'
$ws.Range('E28').Value = '  -3.44%  '

$ws.Range('B39').Value = 'FirstDigitalUSD'
$ws.Range('C39').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.998'
$ws.Range('E39').Value = '  -0.11%  '

$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.130'
$ws.Range('E40').Value = '  -5.27%  '

$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.057.99'
$ws.Range('E41').Value = '  -5.34%  '

$ws.Range('B42').Value = 'ThetaToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.98'
$ws.Range('E42').Value = '  -4.96%  '

